# Applies FlashScore odds corrections for rows 2-5 (Sheet1) per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 11
$ws.Range("Q2").Value = 1.93
$ws.Range("R2").Value = 1.93

# Row 3
$ws.Range("G3").Value = 2.45
$ws.Range("I3").Value = 3.2
$ws.Range("J3").Value = 3.4
$ws.Range("L3").Value = 4
$ws.Range("O3").Value = 1.57
$ws.Range("P3").Value = 2.25
$ws.Range("S3").Value = 1.62
$ws.Range("T3").Value = 2.2
$ws.Range("U3").Value = 2.25
$ws.Range("V3").Value = 1.57
$ws.Range("W3").Value = 6
$ws.Range("X3").Value = 10
$ws.Range("Z3").Value = 23
$ws.Range("AA3").Value = 26
$ws.Range("AC3").Value = 6
$ws.Range("AE3").Value = 21
$ws.Range("AF3").Value = 81
$ws.Range("AI3").Value = 13
$ws.Range("AJ3").Value = 13
$ws.Range("AK3").Value = 34
$ws.Range("AL3").Value = 34
$ws.Range("AN3").Value = 4.33
$ws.Range("AT3").Value = 2.2
$ws.Range("AU3").Value = 9.5
$ws.Range("AW3").Value = 5
$ws.Range("AX3").Value = 21
$ws.Range("AZ3").Value = 67
$ws.Range("BA3").Value = 126
$ws.Range("BB3").Value = 351

# Row 4
$ws.Range("G4").Value = 3.2
$ws.Range("H4").Value = 2.88
$ws.Range("I4").Value = 2.45
$ws.Range("J4").Value = 3.75
$ws.Range("L4").Value = 3.25
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("W4").Value = 8.5
$ws.Range("X4").Value = 15
$ws.Range("Y4").Value = 12
$ws.Range("Z4").Value = 34
$ws.Range("AA4").Value = 29
$ws.Range("AG4").Value = 351
$ws.Range("AH4").Value = 7
$ws.Range("AI4").Value = 11
$ws.Range("AJ4").Value = 10
$ws.Range("AK4").Value = 23
$ws.Range("AN4").Value = 4.75
$ws.Range("AO4").Value = 17
$ws.Range("AP4").Value = 29
$ws.Range("AQ4").Value = 51
$ws.Range("AR4").Value = 81
$ws.Range("AW4").Value = 4.33
$ws.Range("AX4").Value = 15
$ws.Range("AZ4").Value = 51
$ws.Range("BA4").Value = 81
$ws.Range("BB4").Value = 251

# Row 5
$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("Q5").Value = 2.25
$ws.Range("R5").Value = 1.62
$ws.Range("AG5").Value = 900
